# Reset the "Products" and "Sales" sheets back to just their header rows,
# reordering/renaming the header columns to the new schema.

$wb = $excel.ActiveWorkbook

# --- Products sheet (Sheet1) ---
$products = $wb.Worksheets.Item("Products")

# Clear any existing data (header + sample row) before rewriting headers.
$products.Cells.Clear()

$productHeaders = @("id", "name", "description", "price", "cost_price", "category", "stock", "min_stock", "supplier", "sku", "created_date", "last_updated")
for ($i = 0; $i -lt $productHeaders.Length; $i++) {
    $products.Cells.Item(1, $i + 1).Value = $productHeaders[$i]
}

# --- Sales sheet (Sheet2) ---
$sales = $wb.Worksheets.Item("Sales")

$sales.Cells.Clear()

$salesHeaders = @("id", "product_id", "product_name", "quantity", "unit_price", "total_amount", "profit", "customer_name", "payment_method", "sale_date", "cashier", "notes")
for ($i = 0; $i -lt $salesHeaders.Length; $i++) {
    $sales.Cells.Item(1, $i + 1).Value = $salesHeaders[$i]
}
